$d = $word.ActiveDocument
$tbl = $d.Tables(1)
$br = [char]11

# --- Cell (3,5): "Description" column for the requestId row ---
$descCell = $tbl.Cell(3, 5)
$descRange = $descCell.Range
$descRange.Text = "Identifiant unique partagé de la demande de ressource,  généré une seule fois par le système du partenaire qui émet la demande " + $br + "Il est valorisé comme suit lors de sa création : " + $br + "{orgID}.request.{ID unique de la demande dans le système émetteur}" + $br + $br + "OU - uniquement si un ID unique de la demande n'est pas disponible : " + $br + "OrgId émetteur}.request.{senderCaseId}.{numéro d’ordre chronologique}"

# --- Cell (3,6): "Exemple" column for the requestId row ---
$exCell = $tbl.Cell(3, 6)
$exRange = $exCell.Range
$exRange.Text = "fr.health.samu770.request.1249875" + $br + "fr.health.samu690.request.DRFR15690242370035.3"
